$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.596.21'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.470.17'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.89'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.00'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.471.23'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.486'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.122'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.14'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.377'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.073.16'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.08'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.512.21'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000176'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.813.88'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.02'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.21'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.61'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '381.51'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.570'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.621.65'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '72.57'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000110'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.85%  '
$ws.Range('E28').Value = '  -4.31%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.35'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.22'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.12'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.487.99'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.54'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  -3.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.21'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '167.42'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.54'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.82'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0792'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.47'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.808'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('E45').Value = '  -5.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.10'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.34'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.62'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.423.57'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.80'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.879'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.34%  '
